$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:3").Delete()
$ws.Range("A1").Select()
